$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '35.304.02'
$ws.Range("E2").Value = '  +2.01%  '
$ws.Range("D3").Value = '1.845.76'
$ws.Range("E3").Value = '  +1.88%  '
$ws.Range("E4").Value = '  +0.17%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '227.62'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("E6").Value = '  +2.10%  '
$ws.Range("E7").Value = '  +0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '42.63'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +12.53%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.307'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.99%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0687'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.83%  '
$ws.Range("E11").Value = '  +3.68%  '
$ws.Range("D12").Value = '2.114.91'
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '11.65'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.34%  '
$ws.Range("D14").Value = '1.843.22'
$ws.Range("E14").Value = '  +1.52%  '
$ws.Range("E15").Value = '  +6.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.663'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.46%  '
$ws.Range("D17").Value = '35.265.82'
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '69.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.09'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  +1.72%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.14'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.76'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +14.99%  '
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("E24").Value = '  -0.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '170.99'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.92'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.29%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.85'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.87%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.122'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.89%  '
$ws.Range("D29").Value = '3.522.65'
$ws.Range("E29").Value = '  +44.98%  '
$ws.Range("E30").Value = '  +0.24%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.32'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.38%  '
$ws.Range("E32").Value = '  +2.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.02'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0533'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.88'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.82%  '
$ws.Range("E36").Value = '  +3.16%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '88.88'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +9.57%  '
$ws.Range("E38").Value = '  +1.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.02'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.92%  '
$ws.Range("D40").Value = '1.334.95'
$ws.Range("E40").Value = '  -2.34%  '
$ws.Range("E41").Value = '  +3.27%  '
$ws.Range("E42").Value = '  +1.82%  '
$ws.Range("E43").Value = '  +4.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.45'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("E46").Value = '  +1.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0520'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.49%  '
$ws.Range("D48").Value = '2.013.79'
$ws.Range("E48").Value = '  +1.98%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.99'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.73%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '104.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.35%  '
$ws.Range("E51").Value = '  +0.16%  '
